$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Version bump 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date bump
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value now populated
$ws1.Range("B9").Value = "Alvearie Team"

# Former duplicate "Contact" / "No display for ContactDetail" row (row 10) becomes Jurisdiction
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# The second duplicate "Contact" row (row 11) is removed entirely, shifting everything up
$ws1.Rows.Item(11).Delete()

# --- Elements sheet ---
$ws2 = $wb.Worksheets.Item("Elements")

# Root extension element's Short/Definition now reflect the specific extension instead of generic text
$ws2.Range("K2").Value = "Medicare Indicator"
$ws2.Range("L2").Value = "Indicator of Medicare coverage for the member"
